$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: update text of the "outliers" item (shared string kept at same slot,
#     but its wording changes) and grow the row height to fit the longer text. ---
$ws.Range("A17").Value2 = "62. ajouter un test sur les outliers pour les cas excessifs cf. transcripto rainettes 2018 (implémenter sur chaque item, dès la vérification des données, sans prendre en compte la dose-réponse, une détection basée sur le Z-score modifié de Iglewicz, B., & Hoaglin, D. C. (1993). How to detect and handle outliers (Vol. 16). Asq Press. (1115 citations google scholar))"
$ws.Rows.Item(17).RowHeight = 43.2

# --- Rows 19 & 20 swap content: row 20 becomes the new item #65, row 19 becomes
#     the (reworded) Rager item #64. (Order of assignment matters for shared
#     string table ordering, so set A20 before A19.) ---
$ws.Range("A20").Value2 = "65. changer le nom adjpvalue en qvalue dans les sorties pour que ce soit plus clair !!! (pas clair la tricherie dans Rager 2017)"
$ws.Range("A19").Value2 = "64. sur les données in situ de type Rager 2017 avec les plus faibles doses loin du zéro. Si on n'a pas de dose à zéro, demander à l'utilisateur d'indiquer une valeur seuil en dessous de laquelle on considère que la dose est négligeable, du niveau du contrôle et fixer les doses inférieures ou égales à 0"

# --- New rows 21 & 22: two new TODO items, formatted like the other A-column
#     items in this block (copy formatting from A20). ---
$ws.Range("A20").Copy()
$ws.Range("A21").PasteSpecial(-4122)
$ws.Range("A22").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A21").Value2 = "66. éliminer en amont, avant comparason des AICs les pics qui dépassent la gamme du signal observé"
$ws.Range("A22").Value2 = "67. tenter de fixer le f à 0 pour les modèles Gauss-probit et log-Gauss-probit, sans réajuster ou en réajustant à partir des mêmes valeurs, et garder le probit ou log-probit s'il est meilleur en AIC"

$ws.Rows.Item(22).RowHeight = 28.8

# --- Update the view selection to mirror where the user ended up after editing ---
[void]$ws.Range("A26").Select()
